# Update "想去人数" (F column) counts across sheets, as produced by the
# gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(5, 6).Value  = 845
$ws1.Cells.Item(9, 6).Value  = 783
$ws1.Cells.Item(15, 6).Value = 883
$ws1.Cells.Item(16, 6).Value = 9682
$ws1.Cells.Item(17, 6).Value = 594
$ws1.Cells.Item(23, 6).Value = 1740
$ws1.Cells.Item(26, 6).Value = 475
$ws1.Cells.Item(36, 6).Value = 181
$ws1.Cells.Item(39, 6).Value = 90

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(12, 6).Value = 79

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value = 804

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value  = 804
$ws4.Cells.Item(9, 6).Value  = 845
$ws4.Cells.Item(15, 6).Value = 783
$ws4.Cells.Item(20, 6).Value = 883
$ws4.Cells.Item(21, 6).Value = 9682
$ws4.Cells.Item(23, 6).Value = 594
$ws4.Cells.Item(27, 6).Value = 1740
$ws4.Cells.Item(29, 6).Value = 475
$ws4.Cells.Item(31, 6).Value = 79
$ws4.Cells.Item(32, 6).Value = 79
$ws4.Cells.Item(46, 6).Value = 181
